$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update DM_Stat (C) and P_Value (D) columns for rows 2-11, and G4 (Significativo) from "Sí" to "No"

$ws.Range("C2").Value = 0.1044352602029346
$ws.Range("D2").Value = 0.917437593216758

$ws.Range("C3").Value = -0.4885304526651872
$ws.Range("D3").Value = 0.6283104280076994

$ws.Range("C4").Value = -1.82906190066623
$ws.Range("D4").Value = 0.07617278644637326
$ws.Range("G4").Value = "No"

$ws.Range("C5").Value = -0.9331340402400529
$ws.Range("D5").Value = 0.3573305317877868

$ws.Range("C6").Value = -0.8583555709111216
$ws.Range("D6").Value = 0.3967071067502275

$ws.Range("C7").Value = -1.708393380181979
$ws.Range("D7").Value = 0.09667916861347425

$ws.Range("C8").Value = -1.433878655399216
$ws.Range("D8").Value = 0.1607403334695758

$ws.Range("C9").Value = -1.10886489851142
$ws.Range("D9").Value = 0.2752706358671222

$ws.Range("C10").Value = -1.297647867694226
$ws.Range("D10").Value = 0.2031437434944698

$ws.Range("C11").Value = 0.3354258581702169
$ws.Range("D11").Value = 0.739364821765006
